$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For the Price column (D), values such as "560.30" or "0.0000147" look
# numeric. Assigning them directly would let Excel auto-convert the cell to a
# floating point number and silently drop significant trailing/leading zeros.
# Prefixing with a single quote forces Excel to keep the literal text, and
# resetting Style to "Normal" afterwards clears the quote-prefix cell format
# so no extra formatting is left behind on the cell.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '63.443.86'
$ws.Range('E2').Value = '  -0.92%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.720.55'
$ws.Range('E3').Value = '  -1.27%  '

# Row 4
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
Set-TextValue $ws.Range('D5') '560.30'
$ws.Range('E5').Value = '  -2.67%  '

# Row 6
Set-TextValue $ws.Range('D6') '158.04'
$ws.Range('E6').Value = '  -0.51%  '

# Row 7
$ws.Range('E7').Value = '  +0.00%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.594'
$ws.Range('E8').Value = '  -1.46%  '

# Row 9
$ws.Range('E9').Value = '  -2.57%  '

# Row 10
$ws.Range('E10').Value = '  +0.15%  '

# Row 11
Set-TextValue $ws.Range('D11') '5.66'
$ws.Range('E11').Value = '  -3.05%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.374'
$ws.Range('E12').Value = '  -3.39%  '

# Row 13
Set-TextValue $ws.Range('D13') '3.201.03'
$ws.Range('E13').Value = '  -1.31%  '

# Row 14
Set-TextValue $ws.Range('D14') '26.52'
$ws.Range('E14').Value = '  -1.62%  '

# Row 15
Set-TextValue $ws.Range('D15') '63.310.13'
$ws.Range('E15').Value = '  -0.57%  '

# Row 16
Set-TextValue $ws.Range('D16') '0.0000147'
$ws.Range('E16').Value = '  -3.12%  '

# Row 17
Set-TextValue $ws.Range('D17') '2.721.41'
$ws.Range('E17').Value = '  -1.40%  '

# Row 18
Set-TextValue $ws.Range('D18') '12.24'
$ws.Range('E18').Value = '  +0.36%  '

# Row 19
Set-TextValue $ws.Range('D19') '4.67'
$ws.Range('E19').Value = '  -3.87%  '

# Row 20
Set-TextValue $ws.Range('D20') '351.00'
$ws.Range('E20').Value = '  -1.71%  '

# Row 21
$ws.Range('E21').Value = '  -4.18%  '

# Row 22
$ws.Range('E22').Value = '  +0.25%  '

# Row 23
Set-TextValue $ws.Range('D23') '0.514'
$ws.Range('E23').Value = '  -3.44%  '

# Row 24
Set-TextValue $ws.Range('D24') '64.40'
$ws.Range('E24').Value = '  -1.42%  '

# Row 25
Set-TextValue $ws.Range('D25') '0.170'
$ws.Range('E25').Value = '  -0.60%  '

# Row 26
Set-TextValue $ws.Range('D26') '0.999'
$ws.Range('E26').Value = '  +0.01%  '

# Row 27
Set-TextValue $ws.Range('D27') '8.22'
$ws.Range('E27').Value = '  -4.57%  '

# Row 28
Set-TextValue $ws.Range('D28') '0.0₃0883'
$ws.Range('E28').Value = '  -2.72%  '

# Row 29
Set-TextValue $ws.Range('D29') '1.36'
$ws.Range('E29').Value = '  +9.09%  '

# Row 30
$ws.Range('E30').Value = '  +0.14%  '

# Row 31
Set-TextValue $ws.Range('D31') '7.14'
$ws.Range('E31').Value = '  -2.55%  '

# Row 32
Set-TextValue $ws.Range('D32') '165.89'
$ws.Range('E32').Value = '  -2.20%  '

# Row 33
$ws.Range('E33').Value = '  +0.27%  '

# Row 34
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D34') '19.88'
$ws.Range('E34').Value = '  -1.57%  '

# Row 35
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D35') '0.998'
$ws.Range('E35').Value = '  -0.04%  '

# Row 36
Set-TextValue $ws.Range('D36') '4.83'
$ws.Range('E36').Value = '  -2.61%  '

# Row 37
Set-TextValue $ws.Range('D37') '1.79'
$ws.Range('E37').Value = '  -0.58%  '

# Row 38
Set-TextValue $ws.Range('D38') '345.06'
$ws.Range('E38').Value = '  -1.26%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.960'
$ws.Range('E39').Value = '  -4.10%  '

# Row 40
Set-TextValue $ws.Range('D40') '6.08'
$ws.Range('E40').Value = '  -3.09%  '

# Row 41
Set-TextValue $ws.Range('D41') '4.05'
$ws.Range('E41').Value = '  -3.39%  '

# Row 42
Set-TextValue $ws.Range('D42') '38.21'
$ws.Range('E42').Value = '  -2.18%  '

# Row 43
$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D43') '21.38'
$ws.Range('E43').Value = '  -2.39%  '

# Row 44
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D44') '20.82'
$ws.Range('E44').Value = '  -3.26%  '

# Row 45
Set-TextValue $ws.Range('D45') '0.0571'
$ws.Range('E45').Value = '  -3.06%  '

# Row 46
Set-TextValue $ws.Range('D46') '0.626'
$ws.Range('E46').Value = '  -1.13%  '

# Row 47
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D47') '132.12'
$ws.Range('E47').Value = '  -3.36%  '

# Row 48
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D48') '0.997'
$ws.Range('E48').Value = '  -0.14%  '

# Row 49
$ws.Range('B49').Value = 'WhiteBITCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range('D49') '11.07'
$ws.Range('E49').Value = '  +0.12%  '

# Row 50
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D50') '0.0985'
$ws.Range('E50').Value = '  -3.27%  '

# Row 51
Set-TextValue $ws.Range('D51') '0.0246'
$ws.Range('E51').Value = '  -3.72%  '
